# Update "Última actualización" timestamp across all schedule sheets,
# and refresh the latest arrival row on the main LP1912 sheet.

$wb = $excel.ActiveWorkbook

$oldTime = "00:21:43"
$newTime = "01:37:45"

# --- Sheet 1: LP1912 (main data sheet) ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A6").Value = $newTime
$ws1.Range("B6").Value = "03:01"
$ws1.Range("C6").Value = "15_ABASTO"
$ws1.Range("D6").Value = 84

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
